# Generate Report for Handback
# The handoff/handback pipeline just produced a new handback for the
# "59df1511-dce7-46f6-9dc0-d61e9dbaecd1" item on BOTH locale sheets, but the
# handback commit wasn't the latest one on the source branch, so besides
# recording the handback file + timestamp we also log a validation error.

$wb = $excel.ActiveWorkbook

$mdName    = "59df1511-dce7-46f6-9dc0-d61e9dbaecd1.md"
$mdUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b150b454b85a23168f3ec2f129dbc4945a83e5ac/e2e/59df1511-dce7-46f6-9dc0-d61e9dbaecd1.md"
$errorText = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/faee704611e0ac710a9945f068d58038200e6500/e2e/59df1511-dce7-46f6-9dc0-d61e9dbaecd1.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b150b454b85a23168f3ec2f129dbc4945a83e5ac/e2e/59df1511-dce7-46f6-9dc0-d61e9dbaecd1.md."

function Update-HandbackRow {
    param([string]$SheetName, [string]$XlfName, [string]$HandbackDateTime)

    $ws = $wb.Worksheets.Item($SheetName)

    # Latest Target File (I7) becomes a hyperlink to the handback markdown,
    # same as the "Source File Name" hyperlink already on A7.
    $ws.Hyperlinks.Add($ws.Range("I7"), $mdUrl, "", "", $mdName) | Out-Null

    # Latest Handback File (J7)
    $ws.Range("J7").Value = $XlfName

    # Latest Handback DateTime (K7)
    $ws.Range("K7").Value = $HandbackDateTime

    # Error Detail (P7)
    $ws.Range("P7").Value = $errorText
}

Update-HandbackRow "zh-cn" "59df1511-dce7-46f6-9dc0-d61e9dbaecd1.5283ad153ae982e830041532e558d1d45dc1f780.zh-cn.xlf" "2016-09-04 05:02:30"
Update-HandbackRow "de-de" "59df1511-dce7-46f6-9dc0-d61e9dbaecd1.5283ad153ae982e830041532e558d1d45dc1f780.de-de.xlf" "2016-09-04 05:02:37"
